$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain, unstyled text values such as
# "68.13" or "1.735.55" (the latter is not a valid number at all).
# Some of the new values we are writing (e.g. "68.00", "7.420") DO look
# like valid numbers, so a naive assignment would make Excel silently
# convert them to numbers and drop the significant trailing zeros.
# To avoid that we temporarily mark those specific cells as Text before
# writing the value, then restore the cell to the default "Normal" style
# so no stray formatting is left behind (matching the original, unstyled
# inline-string cells).
$numericLookingCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D15", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.499.80'
$ws.Range("E2").Value = '  +5.10%  '
$ws.Range("D3").Value = '1.723.33'
$ws.Range("E3").Value = '  +3.98%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '225.97'
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").Value = '0.5368'
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.2683'
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").Value = '0.06612'
$ws.Range("E9").Value = '  +4.05%  '
$ws.Range("D10").Value = '21.67'
$ws.Range("E10").Value = '  +5.21%  '
$ws.Range("D11").Value = '0.07770'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '4.640'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").Value = '1.715.98'
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("D14").Value = '1.961.80'
$ws.Range("E14").Value = '  +4.03%  '
$ws.Range("D15").Value = '0.5879'
$ws.Range("E15").Value = '  +4.62%  '
$ws.Range("D16").Value = '0.0₅8278'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '68.00'
$ws.Range("E17").Value = '  +3.67%  '
$ws.Range("D18").Value = '27.519.15'
$ws.Range("E18").Value = '  +5.17%  '
$ws.Range("D19").Value = '223.79'
$ws.Range("E19").Value = '  +16.08%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").Value = '4.740'
$ws.Range("E21").Value = '  +1.69%  '
$ws.Range("D22").Value = '10.72'
$ws.Range("E22").Value = '  +2.15%  '
$ws.Range("D23").Value = '6.108'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '1.005'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = '148.45'
$ws.Range("E25").Value = '  +1.97%  '
$ws.Range("E26").Value = '  +2.94%  '
$ws.Range("E27").Value = '  +10.74%  '
$ws.Range("D28").Value = '7.420'
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = '16.73'
$ws.Range("E29").Value = '  +4.75%  '
$ws.Range("D30").Value = '0.05578'
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("D31").Value = '1.306'
$ws.Range("E31").Value = '  +2.42%  '
$ws.Range("D32").Value = '3.574'
$ws.Range("E32").Value = '  +2.69%  '
$ws.Range("D33").Value = '3.465'
$ws.Range("E33").Value = '  +2.79%  '
$ws.Range("D34").Value = '1.667'
$ws.Range("E34").Value = '  +6.24%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.452'
$ws.Range("E35").Value = '  +2.03%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '0.9612'
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("D37").Value = '2.815'
$ws.Range("E37").Value = '  +1.31%  '
$ws.Range("D38").Value = '0.5938'
$ws.Range("E38").Value = '  +4.13%  '
$ws.Range("D39").Value = '0.01648'
$ws.Range("E39").Value = '  +3.52%  '
$ws.Range("D40").Value = '5.888'
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("D41").Value = '0.8587'
$ws.Range("E41").Value = '  +3.01%  '
$ws.Range("D42").Value = '1.060.92'
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("D43").Value = '1.005'
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").Value = '101.59'
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").Value = '1.866.59'
$ws.Range("E45").Value = '  +3.90%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₈115'
$ws.Range("E46").Value = '  +9.32%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '59.07'
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.250'
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.4430'
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("B50").Value = 'Frax'
$ws.Range("C50").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D50").Value = '1.004'
$ws.Range("E50").Value = '  +0.51%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05277'
$ws.Range("E51").Value = '  +0.69%  '

foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).Style = "Normal"
}
